# B1-- and B2 PowerPoint.pptx — apply the "new theme" design change.
#
# What actually happened in the authored edit:
#   1. The table on slide 5 was switched to the built-in PowerPoint table
#      style "Medium Style 2 - Accent 1" ({EE97578D-C36B-47FC-8994-E0075C947577}).
#   2. The deck's theme (bound to the one-and-only slide master) was changed
#      from the custom "Integral" / "Red Violet" palette to the standard
#      Office Theme palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#      (PowerPoint re-used the existing theme part for the notes master
#      slot when doing this; the font scheme and format scheme were
#      already identical between the two themes, so only the colour
#      scheme actually changes visually.)

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------------
$slide5 = $p.Slides.Item(5)
$tableShape = $slide5.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{EE97578D-C36B-47FC-8994-E0075C947577}")

# --- 2. Re-colour the presentation theme to the Office Theme palette -----
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
